$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.43424333333333
$ws.Range("H2").Value = 31.30273
$ws.Range("I2").Value = 0.9711091978791583
$ws.Range("J2").Value = 0.9711091978791584
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.219727
$ws.Range("N2").Value = 6.659181
$ws.Range("O2").Value = 0.1284523376175785
$ws.Range("P2").Value = 0.1284523376175785
$ws.Range("Q2").Value = 23.16117165157
$ws.Range("R2").Value = 208.45054486413
$ws.Range("S2").Value = 0.1247412465495095
$ws.Range("T2").Value = 0.1247412465495095
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.43424333333333
$ws.Range("H3").Value = 31.30273
$ws.Range("I3").Value = 0.9711091978791583
$ws.Range("J3").Value = 0.9711091978791584
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.50585933333333
$ws.Range("N3").Value = 34.517578
$ws.Range("O3").Value = 0.6658271614778302
$ws.Range("P3").Value = 0.6658271614778303
$ws.Range("Q3").Value = 120.0549360431044
$ws.Range("R3").Value = 1080.49442438794
$ws.Range("S3").Value = 0.6465908807088925
$ws.Range("T3").Value = 0.6465908807088926
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.43424333333333
$ws.Range("H4").Value = 31.30273
$ws.Range("I4").Value = 0.9711091978791583
$ws.Range("J4").Value = 0.9711091978791584
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01066666666666667
$ws.Range("N4").Value = 0.032
$ws.Range("O4").Value = 0.0006172643158013742
$ws.Range("P4").Value = 0.0006172643158013743
$ws.Range("Q4").Value = 0.1112985955555555
$ws.Range("R4").Value = 1.00168736
$ws.Range("S4").Value = 0.0005994310545972999
$ws.Range("T4").Value = 0.0005994310545973001
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd8"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.43424333333333
$ws.Range("H5").Value = 31.30273
$ws.Range("I5").Value = 0.9711091978791583
$ws.Range("J5").Value = 0.9711091978791584
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.544296666666666
$ws.Range("N5").Value = 10.63289
$ws.Range("O5").Value = 0.2051032365887898
$ws.Range("P5").Value = 0.2051032365887898
$ws.Range("Q5").Value = 36.98205386552222
$ws.Range("R5").Value = 332.8384847897
$ws.Range("S5").Value = 0.1991776395661589
$ws.Range("T5").Value = 0.1991776395661589
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd8"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.310422
$ws.Range("H6").Value = 0.9312659999999999
$ws.Range("I6").Value = 0.02889080212084161
$ws.Range("J6").Value = 0.02889080212084161
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 2.219727
$ws.Range("N6").Value = 6.659181
$ws.Range("O6").Value = 0.1284523376175785
$ws.Range("P6").Value = 0.1284523376175785
$ws.Range("Q6").Value = 0.689052094794
$ws.Range("R6").Value = 6.201468853145999
$ws.Range("S6").Value = 0.003711091068068999
$ws.Range("T6").Value = 0.003711091068068999
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd8"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.310422
$ws.Range("H7").Value = 0.9312659999999999
$ws.Range("I7").Value = 0.02889080212084161
$ws.Range("J7").Value = 0.02889080212084161
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.50585933333333
$ws.Range("N7").Value = 34.517578
$ws.Range("O7").Value = 0.6658271614778302
$ws.Range("P7").Value = 0.6658271614778303
$ws.Range("Q7").Value = 3.571671865972
$ws.Range("R7").Value = 32.14504679374799
$ws.Range("S7").Value = 0.01923628076893764
$ws.Range("T7").Value = 0.01923628076893765
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd8"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.310422
$ws.Range("H8").Value = 0.9312659999999999
$ws.Range("I8").Value = 0.02889080212084161
$ws.Range("J8").Value = 0.02889080212084161
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01066666666666667
$ws.Range("N8").Value = 0.032
$ws.Range("O8").Value = 0.0006172643158013742
$ws.Range("P8").Value = 0.0006172643158013743
$ws.Range("Q8").Value = 0.003311168
$ws.Range("R8").Value = 0.029800512
$ws.Range("S8").Value = 0.00001783326120407419
$ws.Range("T8").Value = 0.00001783326120407419
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd8"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.310422
$ws.Range("H9").Value = 0.9312659999999999
$ws.Range("I9").Value = 0.02889080212084161
$ws.Range("J9").Value = 0.02889080212084161
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.544296666666666
$ws.Range("N9").Value = 10.63289
$ws.Range("O9").Value = 0.2051032365887898
$ws.Range("P9").Value = 0.2051032365887898
$ws.Range("Q9").Value = 1.10022765986
$ws.Range("R9").Value = 9.902048938739998
$ws.Range("S9").Value = 0.005925597022630886
$ws.Range("T9").Value = 0.005925597022630888
